$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34/35 swap: Filecoin <-> InternetComputer(DFINITY) ---
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "41.550.79"
$ws.Range("D3").Value = "2.254.32"
$ws.Range("D5").Value = "'233.87"
$ws.Range("D7").Value = "'64.88"
$ws.Range("D10").Value = "'59.28"
$ws.Range("D11").Value = "'0.0899"
$ws.Range("D13").Value = "2.589.63"
$ws.Range("D14").Value = "'16.24"
$ws.Range("D15").Value = "'22.60"
$ws.Range("D17").Value = "'5.70"
$ws.Range("D18").Value = "2.257.28"
$ws.Range("D19").Value = "41.439.01"
$ws.Range("D20").Value = "'74.07"
$ws.Range("D21").Value = "'0.0" + [char]8323 + "0922"
$ws.Range("D22").Value = "'6.21"
$ws.Range("D23").Value = "'253.14"
$ws.Range("D26").Value = "'2.32"
$ws.Range("D28").Value = "'173.26"
$ws.Range("D29").Value = "'0.145"
$ws.Range("D30").Value = "'20.49"
$ws.Range("D31").Value = "'1.44"
$ws.Range("D32").Value = "'2.80"
$ws.Range("D34").Value = "'5.04"
$ws.Range("D35").Value = "'4.75"
$ws.Range("D36").Value = "'7.19"
$ws.Range("D38").Value = "'3.93"
$ws.Range("D40").Value = "'0.997"
$ws.Range("D41").Value = "'0.000237"
$ws.Range("D42").Value = "'4.88"
$ws.Range("D44").Value = "'8.83"
$ws.Range("D45").Value = "'18.01"
$ws.Range("D46").Value = "'102.37"
$ws.Range("D47").Value = "'1.23"
$ws.Range("D48").Value = "1.516.82"
$ws.Range("D49").Value = "'0.0944"
$ws.Range("D50").Value = "'1.12"
$ws.Range("D51").Value = "'2.81"

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = "  +5.34%  "
$ws.Range("E3").Value = "  +4.53%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  +4.89%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  +4.62%  "
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  +8.36%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("E23").Value = "  +9.74%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("E32").Value = "  +5.81%  "
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("E34").Value = "  +6.38%  "
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  +3.14%  "
$ws.Range("E38").Value = "  +10.03%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  +50.13%  "
$ws.Range("E42").Value = "  +12.74%  "
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("E44").Value = "  +13.22%  "
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").Value = "  -0.68%  "
